$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 3.900430680208489
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 9.295990156953671

# Row 3
$ws.Range("B3").Value = 0.01514828764759746
$ws.Range("C3").Value = 0.04240448674262143
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 1.35982162114495

# Row 4
$ws.Range("B4").Value = 0.01514828764759746
$ws.Range("C4").Value = 0.3127903958511391
$ws.Range("D4").Value = 337.1190423067083
$ws.Range("E4").Value = 8.660232485948974
$ws.Range("G4").Value = 346.107213476156
